{"js": "// Replace the 100 arithmetic-problem answers in the document's table with\n// their updated values, cell by cell, in document order (row-major).\nconst newValues = [\n  [\"93-52=41\", \"19+58=77\", \"9+76=85\", \"96+3=99\", \"56+14=70\"],\n  [\"21+53=74\", \"59-50=9\", \"41-27=14\", \"80+6=86\", \"18+42=60\"],\n  [\"40+2=42\", \"15+47=62\", \"13+46=59\", \"19+5=24\", \"88-49=39\"],\n  [\"78+7=85\", \"59-50=9\", \"66-18=48\", \"30-20=10\", \"51-24=27\"],\n  [\"22+18=40\", \"8+50=58\", \"42+47=89\", \"67-42=25\", \"68+1=69\"],\n  [\"75-67=8\", \"33+14=47\", \"60+27=87\", \"53+13=66\", \"46+43=89\"],\n  [\"69+7=76\", \"37+5=42\", \"71-59=12\", \"43+40=83\", \"41+15=56\"],\n  [\"35-19=16\", \"37-5=32\", \"22-17=5\", \"44+5=49\", \"78-71=7\"],\n  [\"16+74=90\", \"28+5=33\", \"39+59=98\", \"50+33=83\", \"68+6=74\"],\n  [\"38+19=57\", \"12+10=22\", \"9+77=86\", \"12+74=86\", \"28+66=94\"],\n  [\"15+27=42\", \"33+37=70\", \"75+2=77\", \"34-12=22\", \"70-18=52\"],\n  [\"26+41=67\", \"93-30=63\", \"69-20=49\", \"76-75=1\", \"92-73=19\"],\n  [\"65-54=11\", \"27+15=42\", \"13+79=92\", \"82-42=40\", \"81-7=74\"],\n  [\"29-12=17\", \"7+8=15\", \"22+54=76\", \"66-58=8\", \"5+35=40\"],\n  [\"11+77=88\", \"69+25=94\", \"46-26=20\", \"8+21=29\", \"84-0=84\"],\n  [\"7+85=92\", \"67+6=73\", \"4+31=35\", \"25+14=39\", \"96-22=74\"],\n  [\"95-19=76\", \"8+79=87\", \"16-14=2\", \"81-58=23\", \"24+50=74\"],\n  [\"80-31=49\", \"49+36=85\", \"87-87=0\", \"6+59=65\", \"50+45=95\"],\n  [\"76+10=86\", \"64+27=91\", \"91-78=13\", \"64+12=76\", \"30+40=70\"],\n  [\"31-6=25\", \"48+16=64\", \"0-0=0\", \"70-47=23\", \"86-5=81\"],\n];\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let r = 0; r < rows.items.length; r++) {\n  const row = rows.items[r];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (let c = 0; c < cells.items.length; c++) {\n    cells.items[c].value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables(1)\n$values = @(\n  @(\"93-52=41\",\"19+58=77\",\"9+76=85\",\"96+3=99\",\"56+14=70\"),\n  @(\"21+53=74\",\"59-50=9\",\"41-27=14\",\"80+6=86\",\"18+42=60\"),\n  @(\"40+2=42\",\"15+47=62\",\"13+46=59\",\"19+5=24\",\"88-49=39\"),\n  @(\"78+7=85\",\"59-50=9\",\"66-18=48\",\"30-20=10\",\"51-24=27\"),\n  @(\"22+18=40\",\"8+50=58\",\"42+47=89\",\"67-42=25\",\"68+1=69\"),\n  @(\"75-67=8\",\"33+14=47\",\"60+27=87\",\"53+13=66\",\"46+43=89\"),\n  @(\"69+7=76\",\"37+5=42\",\"71-59=12\",\"43+40=83\",\"41+15=56\"),\n  @(\"35-19=16\",\"37-5=32\",\"22-17=5\",\"44+5=49\",\"78-71=7\"),\n  @(\"16+74=90\",\"28+5=33\",\"39+59=98\",\"50+33=83\",\"68+6=74\"),\n  @(\"38+19=57\",\"12+10=22\",\"9+77=86\",\"12+74=86\",\"28+66=94\"),\n  @(\"15+27=42\",\"33+37=70\",\"75+2=77\",\"34-12=22\",\"70-18=52\"),\n  @(\"26+41=67\",\"93-30=63\",\"69-20=49\",\"76-75=1\",\"92-73=19\"),\n  @(\"65-54=11\",\"27+15=42\",\"13+79=92\",\"82-42=40\",\"81-7=74\"),\n  @(\"29-12=17\",\"7+8=15\",\"22+54=76\",\"66-58=8\",\"5+35=40\"),\n  @(\"11+77=88\",\"69+25=94\",\"46-26=20\",\"8+21=29\",\"84-0=84\"),\n  @(\"7+85=92\",\"67+6=73\",\"4+31=35\",\"25+14=39\",\"96-22=74\"),\n  @(\"95-19=76\",\"8+79=87\",\"16-14=2\",\"81-58=23\",\"24+50=74\"),\n  @(\"80-31=49\",\"49+36=85\",\"87-87=0\",\"6+59=65\",\"50+45=95\"),\n  @(\"76+10=86\",\"64+27=91\",\"91-78=13\",\"64+12=76\",\"30+40=70\"),\n  @(\"31-6=25\",\"48+16=64\",\"0-0=0\",\"70-47=23\",\"86-5=81\")\n)\nfor ($r = 0; $r -lt $values.Length; $r++) {\n  for ($c = 0; $c -lt $values[$r].Length; $c++) {\n    $t.Cell($r+1, $c+1).Range.Text = $values[$r][$c]\n  }\n}\n"}
